$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 32: new entry "Fin logique sélection des joueurs" ---
$ws.Range("B31:G31").Copy()
$ws.Range("B32:G32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B32").Value = 43161
$ws.Range("C32").Value = "Fin logique sélection des joueurs"
$ws.Range("D32").Value = "30 min"
$ws.Range("E32").Value = "Les jouers doivent être sélectionnés pour lancer une partie. Doublon impossible"
$ws.Rows(32).RowHeight = 31.5

# --- Row 33: new entry "Prototype de remplissage de grille" ---
$ws.Range("B31:G31").Copy()
$ws.Range("B33:G33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B33").Value = 43161
$ws.Range("C33").Value = "Prototype de remplissage de grille"
$ws.Range("D33").Value = "1h"
$ws.Range("E33").Value = "Un table layout panel a été placé au dessous du dessin du plateau, de ce fait il sera possible de placer les pions en fonctions d'une ligne et d'une colonne."
$ws.Rows(33).RowHeight = 63

# --- Row 34: B34 gets a new style (center/vertical-center + wrap) ---
$ws.Range("B34").WrapText = $true
$ws.Range("B34").HorizontalAlignment = -4108
$ws.Range("B34").VerticalAlignment = -4108

# --- Update the active selection to match the new working cell ---
[void]$ws.Range("E33:G33").Select()

Write-Output "done"
